$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number need to be forced to text
# (matching the source data which stores these as literal strings),
# otherwise Excel auto-converts the numeric-looking text into a real number.
$textForceCells = @("D5", "D6", "D8", "D9", "D10", "D14", "D15", "D16", "D17", "D19", "D23", "D25", "D26", "D27", "D30", "D32", "D37", "D38", "D40", "D41", "D43", "D44", "D46", "D47", "D50")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.019.69"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.684.17"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "216.06"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.253"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "21.45"
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("D10").Value = "0.0621"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.925.16"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "1.704.13"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "0.535"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "66.28"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "8.24"
$ws.Range("E17").Value = "  +5.38%  "
$ws.Range("D18").Value = "27.076.34"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "237.42"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "9.25"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "146.90"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "7.24"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "16.09"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "1.513.13"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "0.589"
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("D38").Value = "0.918"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "68.47"
$ws.Range("E43").Value = "  +3.95%  "
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "1.827.65"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "90.40"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "7.93"
$ws.Range("E50").Value = "  +4.68%  "
$ws.Range("E51").Value = "  -0.03%  "

foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).ClearFormats()
}
